# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values for the dd20ad19... row
# (row 3) on both the "zh-cn" and "de-de" worksheets, giving that row its
# own handback timestamps instead of sharing the ones from row 2.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-19 04:40:25"
$wsZhCn.Range("H3").Value = "2016-03-19 04:40:44"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-19 04:40:28"
$wsDeDe.Range("H3").Value = "2016-03-19 04:40:48"
